$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3761
$ws.Range("L3").Value = 3928
$ws.Range("F4").Value = 1930
$ws.Range("K4").Value = 1774
$ws.Range("L4").Value = 973
$ws.Range("L6").Value = 3404
$ws.Range("F7").Value = 24124
$ws.Range("K7").Value = 27568
$ws.Range("L7").Value = 12297

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 234
$ws.Range("L3").Value = 266
$ws.Range("L4").Value = 53
$ws.Range("L6").Value = 217
$ws.Range("L7").Value = 798

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 92
$ws.Range("L6").Value = 62
$ws.Range("L7").Value = 279

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 160
$ws.Range("L7").Value = 573

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 136
$ws.Range("L3").Value = 145
$ws.Range("K4").Value = 40
$ws.Range("K7").Value = 908
$ws.Range("L7").Value = 447

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 54
$ws.Range("L3").Value = 86
$ws.Range("L7").Value = 204

$ws = $wb.Worksheets.Item("Fuller Park")
$ws.Range("L6").Value = 27
$ws.Range("L7").Value = 61

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 102
$ws.Range("L4").Value = 48
$ws.Range("L6").Value = 99
$ws.Range("L7").Value = 406
$ws.Range("L8").Value = 798
$ws.Range("L9").Value = 75
$ws.Range("L15").Value = 88
$ws.Range("L17").Value = 23
$ws.Range("L19").Value = 347
$ws.Range("L20").Value = 307
$ws.Range("L21").Value = 37
$ws.Range("L23").Value = 132
$ws.Range("L25").Value = 67
$ws.Range("L29").Value = 677
$ws.Range("L30").Value = 61
$ws.Range("L33").Value = 573
$ws.Range("L36").Value = 164
$ws.Range("K37").Value = 908
$ws.Range("L37").Value = 447
$ws.Range("L42").Value = 392
$ws.Range("L43").Value = 93
$ws.Range("L46").Value = 29
$ws.Range("L47").Value = 91
$ws.Range("L48").Value = 167
$ws.Range("L49").Value = 64
$ws.Range("L54").Value = 256
$ws.Range("L60").Value = 74
$ws.Range("L62").Value = 7
$ws.Range("F63").Value = 215
$ws.Range("I63").Value = 264
$ws.Range("K63").Value = 165
$ws.Range("L63").Value = 45
$ws.Range("L67").Value = 439
$ws.Range("L73").Value = 103
$ws.Range("L76").Value = 181
$ws.Range("L77").Value = 81
$ws.Range("L78").Value = 155
$ws.Range("L79").Value = 322
$ws.Range("L82").Value = 20
$ws.Range("L83").Value = 279
$ws.Range("L84").Value = 119
$ws.Range("I85").Value = 1163
$ws.Range("L85").Value = 644
$ws.Range("L86").Value = 90
$ws.Range("L89").Value = 174
$ws.Range("L91").Value = 178
$ws.Range("L94").Value = 152
$ws.Range("L96").Value = 125
$ws.Range("L99").Value = 204
$ws.Range("F101").Value = 24124
$ws.Range("K101").Value = 27568
$ws.Range("L101").Value = 12297

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L3").Value = 168
$ws.Range("L7").Value = 439

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 119

$ws = $wb.Worksheets.Item("Lincoln Park")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 64

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("L3").Value = 63
$ws.Range("L4").Value = 21
$ws.Range("L7").Value = 256

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 209
$ws.Range("L3").Value = 256
$ws.Range("L7").Value = 677

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("L2").Value = 23
$ws.Range("L7").Value = 167

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L3").Value = 106
$ws.Range("L6").Value = 103
$ws.Range("L7").Value = 347

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L2").Value = 43
$ws.Range("L7").Value = 99

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L2").Value = 119
$ws.Range("L6").Value = 112
$ws.Range("L7").Value = 392

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L3").Value = 47
$ws.Range("L4").Value = 16
$ws.Range("L6").Value = 47
$ws.Range("L7").Value = 155

$ws = $wb.Worksheets.Item("Jefferson Park")
$ws.Range("L3").Value = 10
$ws.Range("L7").Value = 29

$ws = $wb.Worksheets.Item("Douglas")
$ws.Range("L3").Value = 52
$ws.Range("L7").Value = 132

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L3").Value = 32
$ws.Range("L6").Value = 34
$ws.Range("L7").Value = 125

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 64
$ws.Range("L4").Value = 12
$ws.Range("L7").Value = 178

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L3").Value = 9
$ws.Range("L7").Value = 37

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("L2").Value = 112
$ws.Range("L3").Value = 116
$ws.Range("L6").Value = 64
$ws.Range("L7").Value = 322

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L6").Value = 81
$ws.Range("L7").Value = 307

$ws = $wb.Worksheets.Item("Burnside")
$ws.Range("L6").Value = 8
$ws.Range("L7").Value = 23

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("L2").Value = 63
$ws.Range("L7").Value = 164

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L3").Value = 124
$ws.Range("L7").Value = 406

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 22
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("East Side")
$ws.Range("L2").Value = 23
$ws.Range("L3").Value = 33
$ws.Range("L7").Value = 67

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("L6").Value = 22
$ws.Range("L7").Value = 91

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("L2").Value = 35
$ws.Range("L3").Value = 25
$ws.Range("L6").Value = 19
$ws.Range("L7").Value = 88

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("L3").Value = 30
$ws.Range("L7").Value = 75

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L6").Value = 28
$ws.Range("L7").Value = 102

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("L2").Value = 51
$ws.Range("L7").Value = 174

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("L2").Value = 13
$ws.Range("L3").Value = 19
$ws.Range("L7").Value = 90

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("L3").Value = 37
$ws.Range("L6").Value = 30

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("L3").Value = 28
$ws.Range("L7").Value = 74

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("L3").Value = 28
$ws.Range("L6").Value = 30
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L2").Value = 187
$ws.Range("L3").Value = 264
$ws.Range("I5").Value = 38
$ws.Range("I7").Value = 1163
$ws.Range("L7").Value = 644

$ws = $wb.Worksheets.Item("Sheffield & DePaul")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 20

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L4").Value = 4
$ws.Range("L7").Value = 81

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("L6").Value = 17
$ws.Range("L7").Value = 48

$ws = $wb.Worksheets.Item("Museum Campus")
$ws.Range("L2").Value = 5
$ws.Range("L7").Value = 7
